$d = $word.ActiveDocument

$pairs = @(
    @("361÷8=45, 1", "394÷8=49, 2"),
    @("252÷7=36, 0", "698÷9=77, 5"),
    @("845÷2=422, 1", "207÷9=23, 0"),
    @("975÷7=139, 2", "124÷5=24, 4"),
    @("111÷4=27, 3", "124÷9=13, 7"),
    @("552÷9=61, 3", "688÷9=76, 4"),
    @("381÷7=54, 3", "192÷2=96, 0"),
    @("828÷8=103, 4", "654÷6=109, 0"),
    @("837÷5=167, 2", "788÷2=394, 0"),
    @("214÷5=42, 4", "647÷6=107, 5"),
    @("973÷3=324, 1", "137÷4=34, 1"),
    @("307÷6=51, 1", "238÷8=29, 6"),
    @("675÷8=84, 3", "663÷5=132, 3"),
    @("414÷9=46, 0", "556÷5=111, 1"),
    @("801÷3=267, 0", "165÷6=27, 3"),
    @("288÷3=96, 0", "478÷7=68, 2"),
    @("858÷8=107, 2", "700÷9=77, 7"),
    @("640÷8=80, 0", "473÷4=118, 1"),
    @("265÷4=66, 1", "793÷2=396, 1"),
    @("985÷8=123, 1", "125÷5=25, 0"),
    @("534÷9=59, 3", "223÷4=55, 3"),
    @("965÷9=107, 2", "915÷7=130, 5"),
    @("124÷2=62, 0", "742÷8=92, 6"),
    @("214÷8=26, 6", "723÷6=120, 3"),
    @("777÷3=259, 0", "196÷4=49, 0")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
